$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the translation table with 3 new rows (284-286), reusing the
# formatting of the last existing data row (283) for each new row.
$ws.Range("A283:C283").Copy()
$ws.Range("A284:C286").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(284, 1).Value = "cs"
$ws.Cells.Item(284, 2).Value = "lab.build.button.create"
$ws.Cells.Item(284, 3).Value = "Nový build"

$ws.Cells.Item(285, 1).Value = "cs"
$ws.Cells.Item(285, 2).Value = "lab.build.button.list"
$ws.Cells.Item(285, 3).Value = "Seznam buildů"

$ws.Cells.Item(286, 1).Value = "cs"
$ws.Cells.Item(286, 3).Value = "Nejnovější buildy"
$ws.Cells.Item(286, 2).Value = "lab.build.latest.title"

$ws.Range("B279").Select()
